$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Writes $value into $cell as plain text, even when $value looks like a
# number (e.g. "4.5", "-11.234") or a date/time string, without ever
# letting Excel's automatic "looks like a number" conversion turn it into
# a numeric cell and without registering a new (unused) cell style. We do
# this by building the text in a scratch cell via a formula that evaluates
# to a string ( ="..." ), then copy/paste-special-values that already-text
# result into the destination cell - a plain paste of a string never gets
# re-parsed as a number, and a formula cell never needs a NumberFormat
# change to hold text.
function Set-TextValue($cell, $value) {
    $escaped = $value -replace '"', '""'
    $scratch = $ws.Range("ZZ100")
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $cell.PasteSpecial(-4163)
    $scratch.Clear()
}

# Row 2: Latitud/Longitud were re-entered swapped relative to the
# previously uploaded file.
Set-TextValue $ws.Cells.Item(2, 7) "9.532"
Set-TextValue $ws.Cells.Item(2, 8) "-11.453"

# Row 3: new earthquake record
Set-TextValue $ws.Cells.Item(3, 1) "10/10/2010 23:45:32"
Set-TextValue $ws.Cells.Item(3, 2) "4.5"
Set-TextValue $ws.Cells.Item(3, 3) "42,3"
Set-TextValue $ws.Cells.Item(3, 4) "Tectonico_Falla_Local"
Set-TextValue $ws.Cells.Item(3, 5) "Limón"
Set-TextValue $ws.Cells.Item(3, 6) "Reportado en todas las provincias"
Set-TextValue $ws.Cells.Item(3, 7) "-11.234"
Set-TextValue $ws.Cells.Item(3, 8) "4.321"

# Row 4: new earthquake record
Set-TextValue $ws.Cells.Item(4, 1) "01/01/2011 01:01:01"
Set-TextValue $ws.Cells.Item(4, 2) "3.4"
Set-TextValue $ws.Cells.Item(4, 3) "24.2"
Set-TextValue $ws.Cells.Item(4, 4) "Choque_Placas"
Set-TextValue $ws.Cells.Item(4, 5) "Mar_Caribe"
Set-TextValue $ws.Cells.Item(4, 6) "Choque de la placa tectonica coco"
Set-TextValue $ws.Cells.Item(4, 7) "11.45"
Set-TextValue $ws.Cells.Item(4, 8) "12.43"
